# Junction_Flooding_182.xlsx edit:
#  - Row 5 values get rounded to 2 decimal places (custom accuracy)
#  - Row 6 (the old, higher-magnitude sample row) is removed entirely
#  - The sheet's used dimension shrinks from A1:AH6 to A1:AH5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Overwrite row 5 (B5:AH5) with the rounded values.
$row5 = @{
    "B5"  = 1.14
    "C5"  = 0.57
    "D5"  = 0.23
    "E5"  = 1.71
    "F5"  = 1.89
    "G5"  = 0.75
    "H5"  = 10.12
    "I5"  = 0.94
    "J5"  = 0.61
    "K5"  = 0.52
    "L5"  = 0.84
    "M5"  = 1.08
    "N5"  = 0.26
    "O5"  = 0.35
    "P5"  = 1.66
    "Q5"  = 0.29
    "R5"  = 0.25
    "S5"  = 0.02
    "T5"  = 4.1
    "U5"  = 2.63
    "V5"  = 0.82
    "W5"  = 1.87
    "X5"  = 1.19
    "Y5"  = 0.07
    "Z5"  = 4.2
    "AA5" = 0.57
    "AB5" = 0.69
    "AC5" = 1.17
    "AD5" = 0.88
    "AE5" = 0.57
    "AF5" = 9.77
    "AG5" = 0.25
    "AH5" = 0.71
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# 2) Remove the old row 6 entirely (shifts dimension back to A1:AH5).
$ws.Rows("6").Delete()
